$wb = $excel.ActiveWorkbook

# Locate the "ODI Batting" sheet so the new sheet can be inserted right after it.
$odiBatting = $wb.Worksheets.Item("ODI Batting")

# Add the new worksheet right after "ODI Batting"
$newSheet = $wb.Worksheets.Add($null, $odiBatting)
$newSheet.Name = "ODI Batting Extra"

# --- Header row ---
# Reuse the exact header formatting (bold, centered, top-aligned, thin border)
# already present on the "ODI Batting" sheet's header row, instead of setting
# Font/Border properties directly (which would create brand-new style/font
# entries in styles.xml that don't exist in the target workbook).
$odiBatting.Range("A1:F1").Copy()
$newSheet.Range("A1:F1").PasteSpecial(-4122)  # xlPasteFormats
$newSheet.Application.CutCopyMode = $false

$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($col = 1; $col -le $headers.Length; $col++) {
    $newSheet.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# --- Data rows ---
# Columns: A=MATCH_CODE, B=BATTING_POSITION, C=NUM_4, D=NUM_6, E=PERCENT_RUNS_OF_TOTAL, F=MAN_OF_MATCH
# Every column except BATTING_POSITION (B, numeric where present) is stored as text.
$data = @(
    @("3946", "", "", "", "", "NO"),
    @("3948", "", "", "", "", "NO"),
    @("3949", 3, "4", "1", "22.66%", "NO"),
    @("4698", 3, "0", "0", "1.11%", "NO"),
    @("4699", 3, "1", "0", "5.85%", "NO"),
    @("4700", 3, "0", "0", "", "NO")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowIndex = $i + 2
    $rowData = $data[$i]

    # A: MATCH_CODE - text
    $a = $newSheet.Cells.Item($rowIndex, 1)
    $a.Value = "'" + $rowData[0]

    # B: BATTING_POSITION - numeric when present, blank otherwise
    $b = $newSheet.Cells.Item($rowIndex, 2)
    if ($rowData[1] -eq "") {
        $b.Value = ""
    } else {
        $b.Value = $rowData[1]
    }

    # C: NUM_4 - text
    $c = $newSheet.Cells.Item($rowIndex, 3)
    $c.Value = "'" + $rowData[2]

    # D: NUM_6 - text
    $d = $newSheet.Cells.Item($rowIndex, 4)
    $d.Value = "'" + $rowData[3]

    # E: PERCENT_RUNS_OF_TOTAL - text
    $e = $newSheet.Cells.Item($rowIndex, 5)
    if ($rowData[4] -eq "") {
        $e.Value = ""
    } else {
        $e.Value = "'" + $rowData[4]
    }

    # F: MAN_OF_MATCH - text
    $f = $newSheet.Cells.Item($rowIndex, 6)
    $f.Value = $rowData[5]
}

$wb.Worksheets.Item("Player Info").Activate() | Out-Null
$wb.Worksheets.Item("Player Info").Range("A1").Select() | Out-Null
